# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 286
    6  = 589
    11 = 4483
    16 = 129
    20 = 3340
    21 = 77
    22 = 510
    26 = 94
    32 = 638
    33 = 1994
    34 = 364
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
